# Updated AR and SH
$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("covid19_cases_switzerland")
$wsFatal = $wb.Worksheets.Item("covid19_fatalities_switzerland")

# --- Sheet1 (covid19_cases_switzerland) ---

# Row 21: VD (X21) corrected from 2351 to 2215
$wsCases.Range("X21").Value = 2215

# Row 22: add SZ (T22) value, and drop T21 from the AB22 catch-up formula
$wsCases.Range("T22").Value = 99
$wsCases.Range("AB22").Formula = "=SUM(B22:AA22)+X21+S21+P21+C21"

# Row 23: new day (2020-03-27)
$wsCases.Range("A23").Value = 43917
$wsCases.Range("A23").NumberFormat = "yyyy\-mm\-dd;@"
$wsCases.Range("D23").Value = 43
$wsCases.Range("R23").Value = 36
$wsCases.Range("Z23").Value = 94
$wsCases.Range("AB23").Formula = "=AB22+Z23-Z22+R23-R22+D23-D22"

$null = $wsCases.Range("A23").Select()

# --- Sheet2 (covid19_fatalities_switzerland) ---

# Row 22: SG (Q22) corrected from 1 to 2
$wsFatal.Range("Q22").Value = 2

# Row 23: new day (2020-03-27)
$wsFatal.Range("A23").Value = 43917
$wsFatal.Range("A23").NumberFormat = "yyyy\-mm\-dd;@"
$wsFatal.Range("B23").Value = 2
$wsFatal.Range("D23").Value = 2
$wsFatal.Range("E23").Value = 7
$wsFatal.Range("F23").Value = 5
$wsFatal.Range("G23").Value = 12
$wsFatal.Range("H23").Value = 11
$wsFatal.Range("I23").Value = 21
$wsFatal.Range("K23").Value = 9
$wsFatal.Range("M23").Value = 3
$wsFatal.Range("N23").Value = 5
$wsFatal.Range("Q23").Value = 2
$wsFatal.Range("S23").Value = 1
$wsFatal.Range("T23").Value = 1
$wsFatal.Range("U23").Value = 1
$wsFatal.Range("V23").Value = 67
$wsFatal.Range("X23").Value = 21
$wsFatal.Range("Y23").Value = 15
$wsFatal.Range("AA23").Value = 9
$wsFatal.Range("AB23").Formula = "=SUM(B23:AA23)"
$wsFatal.Range("AB23").NumberFormat = "0"

$null = $wsFatal.Range("Y25").Select()
